$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = 1
$ws.Range("F9").Value = -8
$ws.Range("F17").Value = 0
$ws.Range("F18").Value = 2
$ws.Range("F22").Value = 12
$ws.Range("F24").Value = 0
$ws.Range("F26").Value = 0
$ws.Range("F27").Value = -1
$ws.Range("F30").Value = -4
$ws.Range("F31").Value = -6
$ws.Range("F32").Value = 14
$ws.Range("F36").Value = -4
$ws.Range("F37").Value = -1
$ws.Range("F42").Value = -3
$ws.Range("F43").Value = 2
$ws.Range("F49").Value = -1
$ws.Range("F51").Value = -4
$ws.Range("F52").Value = -1
$ws.Range("F55").Value = 0
$ws.Range("F59").Value = -3
$ws.Range("F60").Value = -4
$ws.Range("F63").Value = -3
$ws.Range("F66").Value = -9
$ws.Range("F67").Value = -1
$ws.Range("F68").Value = -6
$ws.Range("F71").Value = -3
$ws.Range("F75").Value = -6
$ws.Range("F76").Value = -8
$ws.Range("F79").Value = 0
$ws.Range("F83").Value = 5
